# Add "-See Portfolio at portfolio.debabrata.xyz" (with a live hyperlink)
# right after the "Key Projects" headline, matching the target diff.

$d = $word.ActiveDocument

# Locate the "Key Projects" heading and collapse to a point right after it.
$headline = $d.Content
$found = $headline.Find.Execute("Key Projects")
if (-not $found) {
    throw "Could not find 'Key Projects' headline"
}
$insertionPoint = $headline.Duplicate
$insertionPoint.Collapse(0)

# Insert the plain-text lead-in. It inherits the run formatting already at
# that location (sz/szCs 36), so it comes out matching the heading's size.
$insertionPoint.InsertAfter([char]0x2014 + "See Portfolio at portfolio.debabrata.xyz")

# Select just the freshly-typed domain text and turn it into a hyperlink.
$linkRange = $d.Content
$linkRange.Find.Execute("portfolio.debabrata.xyz")
$hyperlink = $d.Hyperlinks.Add($linkRange, "http://portfolio.debabrata.xyz/", "", "", "portfolio.debabrata.xyz")

# Keep the hyperlink text visually consistent with the rest of the headline.
$hyperlink.Range.Font.Size = 18

Write-Output "Inserted portfolio hyperlink after Key Projects headline"
